$d = $word.ActiveDocument

# The last paragraph currently ends with "Cf Algo, fibonacci, factorielles."
# and carries the _GoBack bookmark at its very end. In the target revision
# that bookmark moves into the middle of the new "Instanciation" bullet, so
# drop it from here; it gets re-created (baked into the raw XML) below.
$d.Bookmarks("_GoBack").Delete()

function Add-RawParagraph {
    param([string]$xml)

    $last = $d.Paragraphs.Last
    [void]$last.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Style = "Normal"
    $nr = $newPara.Range
    $nr.Collapse(1)
    [void]$nr.InsertXML($xml)
}

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Plain page break paragraph.
Add-RawParagraph "<w:p $w><w:r><w:br w:type=`"page`"/></w:r></w:p>"

# Centered, bold + underlined section title.
Add-RawParagraph "<w:p $w><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Classes et Objets</w:t></w:r></w:p>"

# Bold lead-in line.
Add-RawParagraph "<w:p $w><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Classe :</w:t></w:r></w:p>"

# Bulleted definition list (list style / numbering already defined in the doc).
Add-RawParagraph "<w:p $w><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Entité qui </w:t></w:r><w:r><w:t>définit</w:t></w:r><w:r><w:t xml:space=`"preserve`"> un nouveau type = </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>déclaration</w:t></w:r></w:p>"

Add-RawParagraph "<w:p $w><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>C'est un modèle, qui sert ensuite à créer des trucs…</w:t></w:r><w:r><w:t xml:space=`"preserve`"> = </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>affectation</w:t></w:r></w:p>"

Add-RawParagraph "<w:p $w><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Déclaration : </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Type var ;</w:t></w:r></w:p>"

Add-RawParagraph "<w:p $w><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Instanciation</w:t></w:r><w:r><w:t xml:space=`"preserve`"> (créer en mémoire)</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:r><w:t xml:space=`"preserve`"> : </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>var = new Type() ;</w:t></w:r></w:p>"

Add-RawParagraph "<w:p $w><w:pPr><w:pStyle w:val=`"Paragraphedeliste`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t xml:space=`"preserve`">Affectation : </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>var.champs = … ;</w:t></w:r></w:p>"

Write-Output "done"
